$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.831.40"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.458.22"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.51"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.15"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +12.37%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.458.58"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.22"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("E12").Value = "  +3.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.051.62"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.134"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.14"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.877.41"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.465.16"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.36"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.13"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.82%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.67%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.50"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +10.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.63"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.16"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.63%  "
$ws.Range("E35").Value = "  +11.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.62"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0777"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.967.97"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.44"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.57"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.59"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0318"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.775"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.69"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +10.45%  "
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("E48").Value = "  +8.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "309.86"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.63"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.866"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.59%  "
